# Update NATMI LR-pair edge-weight statistics for Ucn2-Il10rb with
# freshly recomputed TPM-derived values (new TPM input -> rerun of the
# scoring pipeline). Only the numeric score columns change; the
# identifying columns (Sending/Target cluster, Ligand/Receptor symbol)
# and the receptor-expressing-cell columns (K, L) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value2 = 2
$ws.Range("F2").Value2 = 0.6666666666666666
$ws.Range("G2").Value2 = 0.05692066666666667
$ws.Range("H2").Value2 = 0.170762
$ws.Range("I2").Value2 = 0.07235837399807114
$ws.Range("J2").Value2 = 0.07235837399807114
$ws.Range("M2").Value2 = 15.22275333333334
$ws.Range("N2").Value2 = 45.66826
$ws.Range("O2").Value2 = 0.5308207094915162
$ws.Range("P2").Value2 = 0.530820709491516
$ws.Range("Q2").Value2 = 0.8664892682355557
$ws.Range("R2").Value2 = 7.798403414120001
$ws.Range("S2").Value2 = 0.0384093234233086
$ws.Range("T2").Value2 = 0.03840932342330859

# Row 3
$ws.Range("E3").Value2 = 2
$ws.Range("F3").Value2 = 0.6666666666666666
$ws.Range("G3").Value2 = 0.05692066666666667
$ws.Range("H3").Value2 = 0.170762
$ws.Range("I3").Value2 = 0.07235837399807114
$ws.Range("J3").Value2 = 0.07235837399807114
$ws.Range("O3").Value2 = 0.3274110996898236
$ws.Range("P3").Value2 = 0.3274110996898236
$ws.Range("Q3").Value2 = 0.5344520270397778
$ws.Range("R3").Value2 = 4.810068243358
$ws.Range("S3").Value2 = 0.02369093480247601
$ws.Range("T3").Value2 = 0.02369093480247601

# Row 4
$ws.Range("E4").Value2 = 2
$ws.Range("F4").Value2 = 0.6666666666666666
$ws.Range("G4").Value2 = 0.05692066666666667
$ws.Range("H4").Value2 = 0.170762
$ws.Range("I4").Value2 = 0.07235837399807114
$ws.Range("J4").Value2 = 0.07235837399807114
$ws.Range("M4").Value2 = 4.065595333333333
$ws.Range("N4").Value2 = 12.196786
$ws.Range("O4").Value2 = 0.1417681908186603
$ws.Range("P4").Value2 = 0.1417681908186603
$ws.Range("Q4").Value2 = 0.2314163967702222
$ws.Range("R4").Value2 = 2.082747570932
$ws.Range("S4").Value2 = 0.01025811577228654
$ws.Range("T4").Value2 = 0.01025811577228654

# Row 5
$ws.Range("I5").Value2 = 0.6056548703615503
$ws.Range("J5").Value2 = 0.6056548703615503
$ws.Range("M5").Value2 = 15.22275333333334
$ws.Range("N5").Value2 = 45.66826
$ws.Range("O5").Value2 = 0.5308207094915162
$ws.Range("P5").Value2 = 0.530820709491516
$ws.Range("Q5").Value2 = 7.252698152626667
$ws.Range("R5").Value2 = 65.27428337364
$ws.Range("S5").Value2 = 0.3214941479923104
$ws.Range("T5").Value2 = 0.3214941479923103

# Row 6
$ws.Range("I6").Value2 = 0.6056548703615503
$ws.Range("J6").Value2 = 0.6056548703615503
$ws.Range("O6").Value2 = 0.3274110996898236
$ws.Range("P6").Value2 = 0.3274110996898236
$ws.Range("S6").Value2 = 0.1982981271375727
$ws.Range("T6").Value2 = 0.1982981271375727

# Row 7
$ws.Range("I7").Value2 = 0.6056548703615503
$ws.Range("J7").Value2 = 0.6056548703615503
$ws.Range("M7").Value2 = 4.065595333333333
$ws.Range("N7").Value2 = 12.196786
$ws.Range("O7").Value2 = 0.1417681908186603
$ws.Range("P7").Value2 = 0.1417681908186603
$ws.Range("Q7").Value2 = 1.937004109422667
$ws.Range("R7").Value2 = 17.433036984804
$ws.Range("S7").Value2 = 0.08586259523166721
$ws.Range("T7").Value2 = 0.08586259523166721

# Row 8
$ws.Range("G8").Value2 = 0.2532906666666667
$ws.Range("H8").Value2 = 0.7598720000000001
$ws.Range("I8").Value2 = 0.3219867556403786
$ws.Range("J8").Value2 = 0.3219867556403785
$ws.Range("M8").Value2 = 15.22275333333334
$ws.Range("N8").Value2 = 45.66826
$ws.Range("O8").Value2 = 0.5308207094915162
$ws.Range("P8").Value2 = 0.530820709491516
$ws.Range("Q8").Value2 = 3.855781340302224
$ws.Range("R8").Value2 = 34.70203206272001
$ws.Range("S8").Value2 = 0.1709172380758972
$ws.Range("T8").Value2 = 0.1709172380758971

# Row 9
$ws.Range("G9").Value2 = 0.2532906666666667
$ws.Range("H9").Value2 = 0.7598720000000001
$ws.Range("I9").Value2 = 0.3219867556403786
$ws.Range("J9").Value2 = 0.3219867556403785
$ws.Range("O9").Value2 = 0.3274110996898236
$ws.Range("P9").Value2 = 0.3274110996898236
$ws.Range("Q9").Value2 = 2.378252366983112
$ws.Range("R9").Value2 = 21.404271302848
$ws.Range("S9").Value2 = 0.1054220377497749
$ws.Range("T9").Value2 = 0.1054220377497748

# Row 10
$ws.Range("G10").Value2 = 0.2532906666666667
$ws.Range("H10").Value2 = 0.7598720000000001
$ws.Range("I10").Value2 = 0.3219867556403786
$ws.Range("J10").Value2 = 0.3219867556403785
$ws.Range("M10").Value2 = 4.065595333333333
$ws.Range("N10").Value2 = 12.196786
$ws.Range("O10").Value2 = 0.1417681908186603
$ws.Range("P10").Value2 = 0.1417681908186603
$ws.Range("Q10").Value2 = 1.029777352376889
$ws.Range("R10").Value2 = 9.267996171392001
$ws.Range("S10").Value2 = 0.04564747981470653
$ws.Range("T10").Value2 = 0.04564747981470652
